$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "generated at" timestamp in C4 ---
$ws.Range("C4").Value = "2023-12-16T13:11:47"

# --- 2) Append a whole new "region block" worth of rows (156-192) ---
# The document is a repeating template: 4-row data/gap blocks, occasionally
# interrupted by a region-header block (title + column headers). We clone
# those existing patterns from earlier in the sheet instead of re-typing
# every style/merge by hand, then overwrite the bits of data that differ.
# (Row heights are (re)applied at the very end: writing cell values on a
# wrap-text styled row triggers Excel's autofit and clobbers any height we
# set beforehand.)

# 2a) Two 4-row data blocks (156-159, 160-163), cloned from the 11-14 template
$ws.Range("A11:S14").Copy($ws.Range("A156"))
$ws.Range("A11:S14").Copy($ws.Range("A160"))

# 2b) Region-header block (163-168: trailing gap, gap, title, gap, headers,
#     gap), cloned from the 82-87 template
$ws.Range("A82:S87").Copy($ws.Range("A163"))

# 2c) Five more 4-row data blocks (169-172 .. 189-192), cloned from 11-14
$ws.Range("A11:S14").Copy($ws.Range("A169"))
$ws.Range("A11:S14").Copy($ws.Range("A173"))
$ws.Range("A11:S14").Copy($ws.Range("A177"))
$ws.Range("A11:S14").Copy($ws.Range("A181"))
$ws.Range("A11:S14").Copy($ws.Range("A185"))
$ws.Range("A11:S14").Copy($ws.Range("A189"))

# --- 3) Fill in the real row data (overwrite the cloned template values) ---

function Set-DataRow($row, $date, $d, $e, $g, $j, $k, $m) {
    $ws.Range("B$row").Value = $date
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("G$row").Value = $g
    $ws.Range("J$row").Value = $j
    $ws.Range("K$row").Value = $k
    $ws.Range("M$row").Value = $m
}

function Set-OpqRow($row, $o, $p, $q) {
    $ws.Range("O$row").Value = $o
    $ws.Range("P$row").Value = $p
    $ws.Range("Q$row").Value = $q
}

Set-DataRow 156 "2023-11-06T15:22:13" "231" "myTest" "ref1" "mat1" "45000" "another"
Set-OpqRow  157 "55" "321654" "58"

Set-DataRow 160 "2023-11-06T15:34:13" "2222" "myTest" "ref1" "mat2" "36000" "pers"
Set-OpqRow  161 "56" "980980" "61"

$ws.Range("N165").Value = "reg2"

Set-DataRow 169 "2023-11-11T14:32:28" "w" "myTest" "ref1" "mat2" "30000" "pers"
Set-OpqRow  170 "71" "980980" "91"

Set-DataRow 173 "2023-11-11T15:44:23" "w" "myTest" "ref1" "mat2" "6000" "pers"
Set-OpqRow  174 "71" "980980" "96"

Set-DataRow 177 "2023-11-11T16:04:50" "56" "myTest" "ref1" "mat2" "36000" "pers"
Set-OpqRow  178 "72" "980980" "102"

Set-DataRow 181 "2023-11-13T14:49:18" "trt" "myTest" "ref1" "mat2" "36000" "pers"
Set-OpqRow  182 "70" "980980" "103"

Set-DataRow 185 "2023-11-13T14:51:15" "tttt" "myTest" "ref1" "mat2" "45000" "another"
Set-OpqRow  186 "68" "321654" "104"

Set-DataRow 189 "2023-11-13T20:01:16" "gt" "myTest" "ref1" "mat2" "45000" "another"
Set-OpqRow  190 "73" "321654" "105"

# --- 4) Row heights, applied last (see note above) ---
$ws.Rows.Item(155).RowHeight = 5

$ws.Rows.Item(156).RowHeight = 1
$ws.Rows.Item(157).RowHeight = 29
$ws.Rows.Item(158).RowHeight = 1
$ws.Rows.Item(159).RowHeight = 5

$ws.Rows.Item(160).RowHeight = 1
$ws.Rows.Item(161).RowHeight = 29
$ws.Rows.Item(162).RowHeight = 1
$ws.Rows.Item(163).RowHeight = 88

$ws.Rows.Item(164).RowHeight = 20
$ws.Rows.Item(165).RowHeight = 27
$ws.Rows.Item(166).RowHeight = 3
$ws.Rows.Item(167).RowHeight = 20
$ws.Rows.Item(168).RowHeight = 5

$ws.Rows.Item(169).RowHeight = 1
$ws.Rows.Item(170).RowHeight = 29
$ws.Rows.Item(171).RowHeight = 1
$ws.Rows.Item(172).RowHeight = 5

$ws.Rows.Item(173).RowHeight = 1
$ws.Rows.Item(174).RowHeight = 29
$ws.Rows.Item(175).RowHeight = 1
$ws.Rows.Item(176).RowHeight = 5

$ws.Rows.Item(177).RowHeight = 1
$ws.Rows.Item(178).RowHeight = 29
$ws.Rows.Item(179).RowHeight = 1
$ws.Rows.Item(180).RowHeight = 5

$ws.Rows.Item(181).RowHeight = 1
$ws.Rows.Item(182).RowHeight = 29
$ws.Rows.Item(183).RowHeight = 1
$ws.Rows.Item(184).RowHeight = 5

$ws.Rows.Item(185).RowHeight = 1
$ws.Rows.Item(186).RowHeight = 29
$ws.Rows.Item(187).RowHeight = 1
$ws.Rows.Item(188).RowHeight = 5

$ws.Rows.Item(189).RowHeight = 1
$ws.Rows.Item(190).RowHeight = 29
$ws.Rows.Item(191).RowHeight = 1
$ws.Rows.Item(192).RowHeight = 556

Write-Host "Edit complete"
